$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "Tristan"
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "394"
$ws.Range("C2").Value = "Manager"
$ws.Range("D2").Value = "Regular"

# Row 4 (row 3 intentionally left blank)
$ws.Range("A4").Value = "dwfs"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "123"
$ws.Range("C4").Value = "sdfsf"
$ws.Range("D4").Value = "sdfdsf"

# Row 5
$ws.Range("A5").Value = "onad"
$ws.Range("B5").NumberFormat = "@"
$ws.Range("B5").Value = "5345"
$ws.Range("C5").Value = "kdsj"
$ws.Range("D5").Value = "sdf"

$ws.Range("A1").Select()
